$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => D (price) and E (volume/%) new values.
#   D = $null          -> D column unchanged for that row
#   DForceText = $true -> the new D value parses as a plain number (e.g. "584.60",
#                          "0.630", "0.0000125") so Excel would otherwise coerce it
#                          to a numeric value and normalize its textual form
#                          (dropping trailing zeros / using scientific notation).
#                          Pre-formatting the cell as Text ("@") keeps the exact
#                          original text, matching the source data which stores
#                          these prices as literal strings.
$updates = @(
    @{ Row = 2; D = "66.445.83"; DForceText = $false; E = "  -0.27%  " },
    @{ Row = 3; D = "3.463.77"; DForceText = $false; E = "  -1.31%  " },
    @{ Row = 4; D = $null; DForceText = $false; E = "  +0.06%  " },
    @{ Row = 5; D = "584.60"; DForceText = $true; E = "  +0.11%  " },
    @{ Row = 6; D = "177.83"; DForceText = $false; E = "  +1.37%  " },
    @{ Row = 7; D = "0.630"; DForceText = $true; E = "  +5.80%  " },
    @{ Row = 8; D = $null; DForceText = $false; E = "  -0.01%  " },
    @{ Row = 9; D = "3.461.57"; DForceText = $false; E = "  -1.33%  " },
    @{ Row = 10; D = $null; DForceText = $false; E = "  -0.84%  " },
    @{ Row = 11; D = $null; DForceText = $false; E = "  +0.82%  " },
    @{ Row = 12; D = $null; DForceText = $false; E = "  -1.08%  " },
    @{ Row = 13; D = "4.064.01"; DForceText = $false; E = "  -1.14%  " },
    @{ Row = 14; D = $null; DForceText = $false; E = "  +1.44%  " },
    @{ Row = 15; D = "30.17"; DForceText = $false; E = "  -0.79%  " },
    @{ Row = 16; D = "66.341.15"; DForceText = $false; E = "  -0.43%  " },
    @{ Row = 17; D = $null; DForceText = $false; E = "  -0.51%  " },
    @{ Row = 18; D = "3.465.01"; DForceText = $false; E = "  -1.28%  " },
    @{ Row = 19; D = "5.98"; DForceText = $false; E = "  -1.31%  " },
    @{ Row = 20; D = $null; DForceText = $false; E = "  -0.88%  " },
    @{ Row = 21; D = "372.13"; DForceText = $false; E = "  -2.45%  " },
    @{ Row = 22; D = "7.68"; DForceText = $false; E = "  -2.52%  " },
    @{ Row = 23; D = "73.38"; DForceText = $false; E = "  +1.54%  " },
    @{ Row = 24; D = $null; DForceText = $false; E = "  -0.07%  " },
    @{ Row = 25; D = $null; DForceText = $false; E = "  -1.75%  " },
    @{ Row = 26; D = "0.0000125"; DForceText = $true; E = "  +4.27%  " },
    @{ Row = 27; D = "10.03"; DForceText = $false; E = "  +1.69%  " },
    @{ Row = 28; D = $null; DForceText = $false; E = "  +2.97%  " },
    @{ Row = 29; D = $null; DForceText = $false; E = "  -0.12%  " },
    @{ Row = 30; D = "5.96"; DForceText = $false; E = "  +0.86%  " },
    @{ Row = 31; D = $null; DForceText = $false; E = "  -0.82%  " },
    @{ Row = 32; D = "23.73"; DForceText = $false; E = "  -3.42%  " },
    @{ Row = 33; D = "0.999"; DForceText = $false; E = "  -0.02%  " },
    @{ Row = 34; D = $null; DForceText = $false; E = "  -2.63%  " },
    @{ Row = 35; D = "1.27"; DForceText = $false; E = "  -5.80%  " },
    @{ Row = 36; D = $null; DForceText = $false; E = "  -0.84%  " },
    @{ Row = 37; D = "160.97"; DForceText = $false; E = "  -0.33%  " },
    @{ Row = 38; D = "0.885"; DForceText = $false; E = "  -0.80%  " },
    @{ Row = 39; D = "28.10"; DForceText = $true; E = "  -5.98%  " },
    @{ Row = 40; D = "1.81"; DForceText = $false; E = "  +1.39%  " },
    @{ Row = 41; D = "2.822.48"; DForceText = $false; E = "  +3.50%  " },
    @{ Row = 42; D = $null; DForceText = $false; E = "  +0.59%  " },
    @{ Row = 43; D = $null; DForceText = $false; E = "  +1.99%  " },
    @{ Row = 44; D = "6.50"; DForceText = $true; E = "  -0.50%  " },
    @{ Row = 45; D = "0.0694"; DForceText = $false; E = "  -1.25%  " },
    @{ Row = 46; D = "25.18"; DForceText = $false; E = "  +0.69%  " },
    @{ Row = 47; D = "343.21"; DForceText = $false; E = "  +5.80%  " },
    @{ Row = 48; D = $null; DForceText = $false; E = "  -1.57%  " },
    @{ Row = 49; D = "0.0292"; DForceText = $false; E = "  -0.15%  " },
    @{ Row = 50; D = $null; DForceText = $false; E = "  +3.10%  " },
    @{ Row = 51; D = $null; DForceText = $false; E = "  -1.30%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.DForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
